# Generate Report for Archive
#
# This script:
#  1. Replaces the "Ready for handoff" status text with "In Translation"
#     on every sheet/cell where it appears.
#  2. Narrows the "Status" / locale result columns from ~17.22 chars
#     to ~13.41 chars on all three sheets.

$wb = $excel.ActiveWorkbook

# --- Update the status text wherever it is used ---
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value()
        # NOTE: keep the literal on the left of -eq so the comparison type
        # is governed by the string, not by whatever native type the cell
        # happens to hold (e.g. a boolean cell would otherwise coerce the
        # right-hand string into $true and falsely "match").
        if ("Ready for handoff" -eq $v) {
            $cell.Value = "In Translation"
        }
    }
}

# --- Narrow the relevant columns on each sheet ---
# Target stored width is ~13.41 characters; ColumnWidth is snapped to
# whole-pixel boundaries on write (same as real Excel), so 12.5 is the
# input that lands on the closest achievable stored width.
$narrowWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E:E").ColumnWidth = $narrowWidth
$overview.Range("F:F").ColumnWidth = $narrowWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C:C").ColumnWidth = $narrowWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C:C").ColumnWidth = $narrowWidth
